$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(588, 646, 647, 679, 682, 687)
$startRow = 23

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("A28").Select()
